$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 180.27777
$ws.Range("I42").Value = 121.6
$ws.Range("J42").Value = 253.625
$ws.Range("K42").Value = 364.8
$ws.Range("L42").Value = 760.875
$ws.Range("M42").Value = -134.8
$ws.Range("N42").Value = -1220.875

$ws.Range("H74").Value = 6391
$ws.Range("I74").Value = 3991.5454
$ws.Range("K74").Value = 3991.5454
$ws.Range("M74").Value = -3055.5454

$ws.Range("H77").Value = 6391
$ws.Range("I77").Value = 3991.5454
$ws.Range("K77").Value = 19957.727
$ws.Range("M77").Value = -15277.727

$ws.Range("H103").Value = 408.4074
$ws.Range("I103").Value = 461.73685
$ws.Range("J103").Value = 281.75
$ws.Range("K103").Value = 1385.21055
$ws.Range("L103").Value = 845.25
$ws.Range("M103").Value = -799.21055
$ws.Range("N103").Value = -2017.25

$ws.Range("H138").Value = 2751.8628
$ws.Range("J138").Value = 3391.3438
$ws.Range("L138").Value = 10174.0314
$ws.Range("N138").Value = -20454.0314

$ws.Range("H141").Value = 11931.617
$ws.Range("I141").Value = 6375
$ws.Range("J141").Value = 44160
$ws.Range("K141").Value = 19125
$ws.Range("L141").Value = 132480
$ws.Range("M141").Value = -13945
$ws.Range("N141").Value = -142840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5997936
$ws.Range("I45").Value = 11067877
$ws.Range("J45").Value = 6187.727
$ws.Range("K45").Value = 11067877
$ws.Range("L45").Value = 6187.727
$ws.Range("M45").Value = -11067500
$ws.Range("N45").Value = -6941.727

$ws.Range("H61").Value = 6454.6665
$ws.Range("I61").Value = 6527.65
$ws.Range("K61").Value = 6527.65
$ws.Range("M61").Value = -6315.65

$ws.Range("H63").Value = 5096
$ws.Range("I63").Value = 1744.5714
$ws.Range("K63").Value = 1744.5714
$ws.Range("M63").Value = -1058.5714

$ws.Range("H66").Value = 5096
$ws.Range("I66").Value = 1744.5714
$ws.Range("K66").Value = 8722.857
$ws.Range("M66").Value = -5290.857

$ws.Range("H122").Value = 510920.44
$ws.Range("I122").Value = 1640.0333
$ws.Range("J122").Value = 1899867
$ws.Range("K122").Value = 4920.0999
$ws.Range("L122").Value = 5699601
$ws.Range("M122").Value = -2470.0999
$ws.Range("N122").Value = -5704501

$ws.Range("H132").Value = 2553.6897
$ws.Range("I132").Value = 1627.0416
$ws.Range("K132").Value = 4881.1248
$ws.Range("M132").Value = -2351.1248

$ws.Range("H136").Value = 6454.6665
$ws.Range("I136").Value = 6527.65
$ws.Range("K136").Value = 19582.95
$ws.Range("M136").Value = -17032.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5558329.5
$ws.Range("I86").Value = 10003164
$ws.Range("J86").Value = 2287.25
$ws.Range("K86").Value = 10003164
$ws.Range("L86").Value = 2287.25
$ws.Range("M86").Value = -10002041
$ws.Range("N86").Value = -4533.25

$ws.Range("H89").Value = 5558329.5
$ws.Range("I89").Value = 10003164
$ws.Range("J89").Value = 2287.25
$ws.Range("K89").Value = 50015820
$ws.Range("L89").Value = 11436.25
$ws.Range("M89").Value = -50010204
$ws.Range("N89").Value = -22668.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1772.2727
$ws.Range("I16").Value = 1320
$ws.Range("J16").Value = 2149.1667
$ws.Range("K16").Value = 1320
$ws.Range("L16").Value = 2149.1667
$ws.Range("M16").Value = -1033
$ws.Range("N16").Value = -2723.1667

$ws.Range("H31").Value = 22642.428
$ws.Range("I31").Value = 2813.348
$ws.Range("K31").Value = 2813.348
$ws.Range("M31").Value = -2518.348

$ws.Range("H34").Value = 22642.428
$ws.Range("I34").Value = 2813.348
$ws.Range("K34").Value = 2813.348
$ws.Range("M34").Value = -2611.348

$ws.Range("H99").Value = 3417.56
$ws.Range("I99").Value = 3292.4119
$ws.Range("J99").Value = 3683.5
$ws.Range("K99").Value = 3292.4119
$ws.Range("L99").Value = 3683.5
$ws.Range("M99").Value = -1794.4119
$ws.Range("N99").Value = -6679.5

$ws.Range("H113").Value = 1772.2727
$ws.Range("I113").Value = 1320
$ws.Range("J113").Value = 2149.1667
$ws.Range("K113").Value = 1320
$ws.Range("L113").Value = 2149.1667
$ws.Range("M113").Value = 850
$ws.Range("N113").Value = -6489.1667

$ws.Range("H122").Value = 3272.8948
$ws.Range("I122").Value = 3288.0557
$ws.Range("K122").Value = 9864.167099999999
$ws.Range("M122").Value = -7414.167099999999

$ws.Range("H126").Value = 3417.56
$ws.Range("I126").Value = 3292.4119
$ws.Range("J126").Value = 3683.5
$ws.Range("K126").Value = 9877.235700000001
$ws.Range("L126").Value = 11050.5
$ws.Range("M126").Value = -7407.235700000001
$ws.Range("N126").Value = -15990.5

$ws.Range("H141").Value = 156549.2
$ws.Range("J141").Value = 156549.2
$ws.Range("L141").Value = 156549.2
$ws.Range("N141").Value = -166909.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 56277.562
$ws.Range("J12").Value = 1071.7
$ws.Range("L12").Value = 3215.1
$ws.Range("N12").Value = -3561.1

$ws.Range("H103").Value = 335.14285
$ws.Range("I103").Value = 307.66666
$ws.Range("K103").Value = 922.9999799999999
$ws.Range("M103").Value = -43.99997999999994

$ws.Range("H132").Value = 971.2353000000001
$ws.Range("J132").Value = 649.5
$ws.Range("L132").Value = 5845.5
$ws.Range("N132").Value = -10905.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4354295
$ws.Range("I102").Value = 5556482.5
$ws.Range("K102").Value = 5556482.5
$ws.Range("M102").Value = -5554860.5

$ws.Range("H107").Value = 655.4
$ws.Range("J107").Value = 508.57144
$ws.Range("L107").Value = 508.57144
$ws.Range("N107").Value = -4348.57144

$ws.Range("H113").Value = 16670517
$ws.Range("I113").Value = 55557056
$ws.Range("J113").Value = 4857.143
$ws.Range("K113").Value = 55557056
$ws.Range("L113").Value = 4857.143
$ws.Range("M113").Value = -55554886
$ws.Range("N113").Value = -9197.143

$ws.Range("H122").Value = 390975.56
$ws.Range("I122").Value = 559484.8
$ws.Range("J122").Value = 5811.5713
$ws.Range("K122").Value = 1678454.4
$ws.Range("L122").Value = 17434.7139
$ws.Range("M122").Value = -1676004.4
$ws.Range("N122").Value = -22334.7139

$ws.Range("H126").Value = 4443772.5
$ws.Range("I126").Value = 2675927.5
$ws.Range("J126").Value = 6948219
$ws.Range("K126").Value = 8027782.5
$ws.Range("L126").Value = 20844657
$ws.Range("M126").Value = -8025312.5
$ws.Range("N126").Value = -20849597

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10948.385
$ws.Range("J7").Value = 11145.091
$ws.Range("L7").Value = 11145.091
$ws.Range("N7").Value = -11369.091

$ws.Range("H40").Value = 2571.8333
$ws.Range("I40").Value = 1811.3043
$ws.Range("J40").Value = 5070.7144
$ws.Range("K40").Value = 1811.3043
$ws.Range("L40").Value = 5070.7144
$ws.Range("M40").Value = -1675.3043
$ws.Range("N40").Value = -5342.7144

$ws.Range("H46").Value = 2562234.8
$ws.Range("I46").Value = 21739404
$ws.Range("K46").Value = 21739404
$ws.Range("M46").Value = -21739216

$ws.Range("H126").Value = 10948.385
$ws.Range("J126").Value = 11145.091
$ws.Range("L126").Value = 33435.273
$ws.Range("N126").Value = -38375.273

$ws.Range("H132").Value = 11204.429
$ws.Range("I132").Value = 11960.777
$ws.Range("J132").Value = 6666.3335
$ws.Range("K132").Value = 35882.331
$ws.Range("L132").Value = 19999.0005
$ws.Range("M132").Value = -33352.331
$ws.Range("N132").Value = -25059.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7279.95
$ws.Range("I62").Value = 4095
$ws.Range("J62").Value = 7361.615
$ws.Range("K62").Value = 4095
$ws.Range("L62").Value = 7361.615
$ws.Range("M62").Value = -3471
$ws.Range("N62").Value = -8609.615

$ws.Range("H65").Value = 7279.95
$ws.Range("I65").Value = 4095
$ws.Range("J65").Value = 7361.615
$ws.Range("K65").Value = 20475
$ws.Range("L65").Value = 36808.075
$ws.Range("M65").Value = -17355
$ws.Range("N65").Value = -43048.075

$ws.Range("H136").Value = 4019.1738
$ws.Range("I136").Value = 3906.4092
$ws.Range("K136").Value = 11719.2276
$ws.Range("M136").Value = -9169.2276
